# Protocols.xlsx update:
#  - The "permission" instructions text is rewritten with clearer wording,
#    and the old "upload a thumbnail" text (that used to share the same
#    shared-string slot layout) is kept as-is but now occupies the other
#    slot. Net effect per sheet: the cell that explained the permission
#    rule gets the new, longer wording (and grows from a 2-line to a
#    3-line row), while the Person sheet's thumbnail-upload cell and
#    permission-rule cell swap which text block they show so the
#    thumbnail text ends up in C9 and the permission text in C15.
#  - The active sheet/tab moves from "Person" to "Picturestory", and the
#    remembered selection on the Person sheet moves from E10 to C15.

$wb = $excel.ActiveWorkbook

$newPermissionText = 'If you have permission to use the material, indicate "yes"; if you intend to request use permission, indicate "request use permission"; otherwise indicate "no". Only materials for which you indicated you have use permission will be made available to end users; other materials will remain invisible for end users.'
$thumbnailText = 'Upload a thumbnail of the image of the person here of max. 400X400 pixels (or about 300kb); this will be displayed alongside the entry in the repository. Only upload an image here if you have use permission for that image; if this is not the case, do not upload anything.'

# Sheets (and the row holding the "Permission" explanation cell in column C)
# that simply get the reworded permission text and a taller row.
$permissionCells = @{
    "Text"         = 24
    "Film"         = 23
    "Image"        = 17
    "Infographic"  = 18
    "Music"        = 19
    "Picturestory" = 21
}

foreach ($sheetName in $permissionCells.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = $permissionCells[$sheetName]
    $ws.Cells.Item($row, 3).Value = $newPermissionText
    $ws.Rows.Item($row).RowHeight = 48
}

# Person sheet: the thumbnail-upload cell (C9) and the permission cell
# (C15) trade places content-wise, with the permission cell also picking
# up the new wording.
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Cells.Item(9, 3).Value = $thumbnailText
$wsPerson.Cells.Item(15, 3).Value = $newPermissionText
$wsPerson.Rows.Item(15).RowHeight = 48

# Remembered selection on the Person sheet moves to C15.
$wsPerson.Range("C15").Select()

# The active tab moves from Person to Picturestory.
$wsPicturestory = $wb.Worksheets.Item("Picturestory")
$wsPicturestory.Activate()
